$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "INTEREST.RATE"
$ws.Range("H1").Value = "INTEND.DATE"
$ws.Range("I1").Value = "CUST.REMARKS:1"
$ws.Range("J1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("K1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("L1").Value = "PRIN.LIQ.ACCT"
$ws.Range("M1").Value = "INT.LIQ.ACCT"
$ws.Range("N1").Value = "CHRG.LIQ.ACCT"
$ws.Range("O1").Value = "FINAL.MATURITY"
$ws.Range("P1").Value = "EXP.DATE"

$ws.Range("F7").Select()
